# Apply edit: rename "username"/"jamsgra" column to "email"/<email>,
# restyle the email cells with the Trebuchet MS font (replacing a stale
# numeric-format style), and move the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------
$ws.Range("A1").Value = "email"
$ws.Range("A2").Value = "jamsgra.doey@gmail.com"
$ws.Range("A3").Value = "jamsgra.doey@gmail.com"

# --- Style edits ------------------------------------------------------
# A2 previously carried a lingering numeric-format style; swap it for a
# plain Trebuchet MS font, then replicate that same style onto A3 via a
# format-only copy/paste so both rows share one consolidated style.
$ws.Range("A2").Font.Name = "Trebuchet MS"
$ws.Range("A2").NumberFormat = "general"

$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0 | Out-Null

# --- Selection ---------------------------------------------------------
$ws.Range("A3").Select() | Out-Null
